# Reinsurer share fix: take attachment point / limit occurrence into account.
# - Rename header B1 on the "sections" sheet from "cession_rate" to "cession_PCT".
# - For rows 7-36, the attachment_point_100 (C) and limit_occurrence_100 (D)
#   values were recorded swapped; swap them back.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sections")

# Rename the header.
$ws.Range("B1").Value = "cession_PCT"

# Swap columns C and D (attachment_point_100 / limit_occurrence_100) for rows 7 through 36.
for ($row = 7; $row -le 36; $row++) {
    $cCell = $ws.Cells.Item($row, 3)
    $dCell = $ws.Cells.Item($row, 4)

    $cValue = $cCell.Value()
    $dValue = $dCell.Value()

    $cCell.Value = $dValue
    $dCell.Value = $cValue
}
